$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)

# Add new rows 103-128 for twelve_bar_breakout signal (Signal Instances sheet)
# Row 103
$ws1.Cells.Item(103, 1).Value = 'twelve_bar_breakout'
$ws1.Cells.Item(103, 2).Value = 832
$ws1.Cells.Item(103, 4).Value = 0.0178407412022352
$ws1.Cells.Item(103, 5).Value = 466.0378516784252
$ws1.Cells.Item(103, 6).Value = $true
$ws1.Cells.Item(103, 7).Value = 947.1698362444597
$ws1.Cells.Item(103, 8).Value = 33
$ws1.Cells.Item(103, 9).Value = -15.09433635924169
$ws1.Cells.Item(103, 10).Value = 0
$ws1.Cells.Item(103, 11).Value = $false
$ws1.Cells.Item(103, 14).Value = 481.1319845660345

# Row 104
$ws1.Cells.Item(104, 1).Value = 'twelve_bar_breakout'
$ws1.Cells.Item(104, 2).Value = 843
$ws1.Cells.Item(104, 4).Value = 0.0299589857459068
$ws1.Cells.Item(104, 5).Value = 371.9099871749187
$ws1.Cells.Item(104, 6).Value = $true
$ws1.Cells.Item(104, 7).Value = 523.5954114627168
$ws1.Cells.Item(104, 8).Value = 22
$ws1.Cells.Item(104, 9).Value = -25.84271550316049
$ws1.Cells.Item(104, 10).Value = 0
$ws1.Cells.Item(104, 11).Value = $false
$ws1.Cells.Item(104, 14).Value = 151.6854242877981

# Row 105
$ws1.Cells.Item(105, 1).Value = 'twelve_bar_breakout'
$ws1.Cells.Item(105, 2).Value = 844
$ws1.Cells.Item(105, 4).Value = 0.036354724317789
$ws1.Cells.Item(105, 5).Value = 291.6666401950965
$ws1.Cells.Item(105, 6).Value = $true
$ws1.Cells.Item(105, 7).Value = 413.8888107063117
$ws1.Cells.Item(105, 8).Value = 21
$ws1.Cells.Item(105, 9).Value = -25.92593370056484
$ws1.Cells.Item(105, 10).Value = 0
$ws1.Cells.Item(105, 11).Value = $false
$ws1.Cells.Item(105, 14).Value = 122.2221705112152

# Row 106
$ws1.Cells.Item(106, 1).Value = 'twelve_bar_breakout'
$ws1.Cells.Item(106, 2).Value = 889
$ws1.Cells.Item(106, 4).Value = 0.1110838800668716
$ws1.Cells.Item(106, 5).Value = 192.4241617349992
$ws1.Cells.Item(106, 6).Value = $true
$ws1.Cells.Item(106, 7).Value = 259.0907725698341
$ws1.Cells.Item(106, 8).Value = 46
$ws1.Cells.Item(106, 9).Value = -7.272748125917099
$ws1.Cells.Item(106, 10).Value = 0
$ws1.Cells.Item(106, 11).Value = $false
$ws1.Cells.Item(106, 14).Value = 66.66661083483493

# Row 107
$ws1.Cells.Item(107, 1).Value = 'twelve_bar_breakout'
$ws1.Cells.Item(107, 2).Value = 893
$ws1.Cells.Item(107, 4).Value = 0.1622497290372848
$ws1.Cells.Item(107, 5).Value = 72.19919005509998
$ws1.Cells.Item(107, 6).Value = $true
$ws1.Cells.Item(107, 7).Value = 145.8506189807024
$ws1.Cells.Item(107, 8).Value = 42
$ws1.Cells.Item(107, 9).Value = -21.36926827332137
$ws1.Cells.Item(107, 10).Value = 3
$ws1.Cells.Item(107, 11).Value = $false
$ws1.Cells.Item(107, 14).Value = 73.65142892560242

# Row 108
$ws1.Cells.Item(108, 1).Value = 'twelve_bar_breakout'
$ws1.Cells.Item(108, 2).Value = 934
$ws1.Cells.Item(108, 4).Value = 0.316757321357727
$ws1.Cells.Item(108, 5).Value = 319.766230999871
$ws1.Cells.Item(108, 6).Value = $true
$ws1.Cells.Item(108, 7).Value = 343.6769599428212
$ws1.Cells.Item(108, 8).Value = 47
$ws1.Cells.Item(108, 9).Value = -16.57811249903101
$ws1.Cells.Item(108, 10).Value = 9
$ws1.Cells.Item(108, 11).Value = $true
$ws1.Cells.Item(108, 12).Value = 50
$ws1.Cells.Item(108, 13).Value = 259.1923762011446
$ws1.Cells.Item(108, 14).Value = 23.91072894295019
$ws1.Cells.Item(108, 15).Value = -60.57385479872636
$ws1.Cells.Item(108, 16).Value = -84.48458374167654

# Row 109
$ws1.Cells.Item(109, 1).Value = 'twelve_bar_breakout'
$ws1.Cells.Item(109, 2).Value = 964
$ws1.Cells.Item(109, 4).Value = 0.5308461785316467
$ws1.Cells.Item(109, 5).Value = 376.8548282692319
$ws1.Cells.Item(109, 6).Value = $true
$ws1.Cells.Item(109, 7).Value = 454.8509942930094
$ws1.Cells.Item(109, 8).Value = 46
$ws1.Cells.Item(109, 9).Value = -1.77551844535123
$ws1.Cells.Item(109, 10).Value = 3
$ws1.Cells.Item(109, 11).Value = $true
$ws1.Cells.Item(109, 12).Value = 20
$ws1.Cells.Item(109, 13).Value = 114.3310426615583
$ws1.Cells.Item(109, 14).Value = 77.9961660237775
$ws1.Cells.Item(109, 15).Value = -262.5237856076736
$ws1.Cells.Item(109, 16).Value = -340.5199516314511

# Row 110
$ws1.Cells.Item(110, 1).Value = 'twelve_bar_breakout'
$ws1.Cells.Item(110, 2).Value = 1004
$ws1.Cells.Item(110, 4).Value = 1.969213962554932
$ws1.Cells.Item(110, 5).Value = -67.52137422177323
$ws1.Cells.Item(110, 6).Value = $false
$ws1.Cells.Item(110, 7).Value = 124.7863026165353
$ws1.Cells.Item(110, 8).Value = 32
$ws1.Cells.Item(110, 9).Value = -68.80341800461765
$ws1.Cells.Item(110, 10).Value = 52
$ws1.Cells.Item(110, 11).Value = $true
$ws1.Cells.Item(110, 12).Value = 42
$ws1.Cells.Item(110, 13).Value = 9.230772583559089
$ws1.Cells.Item(110, 14).Value = 192.3076768383086
$ws1.Cells.Item(110, 15).Value = 76.75214680533232
$ws1.Cells.Item(110, 16).Value = -115.5555300329763

# Row 111
$ws1.Cells.Item(111, 1).Value = 'twelve_bar_breakout'
$ws1.Cells.Item(111, 2).Value = 1029
$ws1.Cells.Item(111, 4).Value = 3.090152025222778
$ws1.Cells.Item(111, 5).Value = -71.24183427606911
$ws1.Cells.Item(111, 6).Value = $false
$ws1.Cells.Item(111, 7).Value = 43.24613225838542
$ws1.Cells.Item(111, 8).Value = 7
$ws1.Cells.Item(111, 9).Value = -87.36383748476962
$ws1.Cells.Item(111, 10).Value = 28
$ws1.Cells.Item(111, 11).Value = $true
$ws1.Cells.Item(111, 12).Value = 17
$ws1.Cells.Item(111, 13).Value = -30.39217463849536
$ws1.Cells.Item(111, 14).Value = 114.4879665344545
$ws1.Cells.Item(111, 15).Value = 40.84965963757375
$ws1.Cells.Item(111, 16).Value = -73.63830689688078

# Row 112
$ws1.Cells.Item(112, 1).Value = 'twelve_bar_breakout'
$ws1.Cells.Item(112, 2).Value = 1070
$ws1.Cells.Item(112, 4).Value = 0.9526280760765076
$ws1.Cells.Item(112, 5).Value = 72.79158817169862
$ws1.Cells.Item(112, 6).Value = $true
$ws1.Cells.Item(112, 7).Value = 96.81981008376262
$ws1.Cells.Item(112, 8).Value = 46
$ws1.Cells.Item(112, 9).Value = -22.96819754822275
$ws1.Cells.Item(112, 10).Value = 4
$ws1.Cells.Item(112, 11).Value = $false
$ws1.Cells.Item(112, 14).Value = 24.028221912064

# Row 113
$ws1.Cells.Item(113, 1).Value = 'twelve_bar_breakout'
$ws1.Cells.Item(113, 2).Value = 1085
$ws1.Cells.Item(113, 4).Value = 1.070444345474243
$ws1.Cells.Item(113, 5).Value = 30.50316174395609
$ws1.Cells.Item(113, 6).Value = $true
$ws1.Cells.Item(113, 7).Value = 75.15723989440222
$ws1.Cells.Item(113, 8).Value = 31
$ws1.Cells.Item(113, 9).Value = -16.9811302139
$ws1.Cells.Item(113, 10).Value = 0
$ws1.Cells.Item(113, 11).Value = $true
$ws1.Cells.Item(113, 12).Value = 50
$ws1.Cells.Item(113, 13).Value = 16.9811631734283
$ws1.Cells.Item(113, 14).Value = 44.65407815044613
$ws1.Cells.Item(113, 15).Value = -13.52199857052779
$ws1.Cells.Item(113, 16).Value = -58.17607672097392

# Row 114
$ws1.Cells.Item(114, 1).Value = 'twelve_bar_breakout'
$ws1.Cells.Item(114, 2).Value = 1115
$ws1.Cells.Item(114, 4).Value = 1.494582891464233
$ws1.Cells.Item(114, 5).Value = 47.29730958478864
$ws1.Cells.Item(114, 6).Value = $true
$ws1.Cells.Item(114, 7).Value = 56.08104313076475
$ws1.Cells.Item(114, 8).Value = 44
$ws1.Cells.Item(114, 9).Value = -26.12611859532427
$ws1.Cells.Item(114, 10).Value = 20
$ws1.Cells.Item(114, 11).Value = $true
$ws1.Cells.Item(114, 12).Value = 20
$ws1.Cells.Item(114, 13).Value = -16.21620629999513
$ws1.Cells.Item(114, 14).Value = 8.783733545976119
$ws1.Cells.Item(114, 15).Value = -63.51351588478377
$ws1.Cells.Item(114, 16).Value = -72.29724943075989

# Row 115
$ws1.Cells.Item(115, 1).Value = 'twelve_bar_breakout'
$ws1.Cells.Item(115, 2).Value = 1147
$ws1.Cells.Item(115, 4).Value = 1.494582891464233
$ws1.Cells.Item(115, 5).Value = 34.40700556013294
$ws1.Cells.Item(115, 6).Value = $true
$ws1.Cells.Item(115, 7).Value = 65.31527764083806
$ws1.Cells.Item(115, 8).Value = 23
$ws1.Cells.Item(115, 9).Value = -3.603606339397119
$ws1.Cells.Item(115, 10).Value = 0
$ws1.Cells.Item(115, 11).Value = $true
$ws1.Cells.Item(115, 12).Value = 30
$ws1.Cells.Item(115, 13).Value = 33.04933132442191
$ws1.Cells.Item(115, 14).Value = 30.90827208070512
$ws1.Cells.Item(115, 15).Value = -1.357674235711031
$ws1.Cells.Item(115, 16).Value = -32.26594631641615

# Row 116
$ws1.Cells.Item(116, 1).Value = 'twelve_bar_breakout'
$ws1.Cells.Item(116, 2).Value = 1225
$ws1.Cells.Item(116, 4).Value = 1.783771514892578
$ws1.Cells.Item(116, 5).Value = -3.96831994953004
$ws1.Cells.Item(116, 6).Value = $false
$ws1.Cells.Item(116, 7).Value = 18.79822858937822
$ws1.Cells.Item(116, 8).Value = 6
$ws1.Cells.Item(116, 9).Value = -45.37469068370795
$ws1.Cells.Item(116, 10).Value = 30
$ws1.Cells.Item(116, 11).Value = $true
$ws1.Cells.Item(116, 12).Value = 11
$ws1.Cells.Item(116, 13).Value = 8.225558965366082
$ws1.Cells.Item(116, 14).Value = 22.76654853890826
$ws1.Cells.Item(116, 15).Value = 12.19387891489612
$ws1.Cells.Item(116, 16).Value = -10.57266962401213

# Row 117
$ws1.Cells.Item(117, 1).Value = 'twelve_bar_breakout'
$ws1.Cells.Item(117, 2).Value = 1231
$ws1.Cells.Item(117, 4).Value = 2.057367563247681
$ws1.Cells.Item(117, 5).Value = -35.78450556104401
$ws1.Cells.Item(117, 6).Value = $false
$ws1.Cells.Item(117, 7).Value = 3.000018063334668
$ws1.Cells.Item(117, 8).Value = 0
$ws1.Cells.Item(117, 9).Value = -52.63895839944887
$ws1.Cells.Item(117, 10).Value = 24
$ws1.Cells.Item(117, 11).Value = $true
$ws1.Cells.Item(117, 12).Value = 5
$ws1.Cells.Item(117, 13).Value = -6.166660389548309
$ws1.Cells.Item(117, 14).Value = 38.78452362437869
$ws1.Cells.Item(117, 15).Value = 29.61784517149571
$ws1.Cells.Item(117, 16).Value = -9.166678452882977

# Row 118
$ws1.Cells.Item(118, 1).Value = 'twelve_bar_breakout'
$ws1.Cells.Item(118, 2).Value = 1269
$ws1.Cells.Item(118, 4).Value = 1.494528293609619
$ws1.Cells.Item(118, 5).Value = 27.30193605940554
$ws1.Cells.Item(118, 6).Value = $true
$ws1.Cells.Item(118, 7).Value = 40.83962992994832
$ws1.Cells.Item(118, 8).Value = 49
$ws1.Cells.Item(118, 9).Value = -33.41068169660229
$ws1.Cells.Item(118, 10).Value = 28
$ws1.Cells.Item(118, 11).Value = $true
$ws1.Cells.Item(118, 12).Value = 16
$ws1.Cells.Item(118, 13).Value = -17.16939411081364
$ws1.Cells.Item(118, 14).Value = 13.53769387054278
$ws1.Cells.Item(118, 15).Value = -44.47133017021919
$ws1.Cells.Item(118, 16).Value = -58.00902404076197

# Row 119
$ws1.Cells.Item(119, 1).Value = 'twelve_bar_breakout'
$ws1.Cells.Item(119, 2).Value = 1310
$ws1.Cells.Item(119, 4).Value = 1.753481864929199
$ws1.Cells.Item(119, 5).Value = -25.26360244069183
$ws1.Cells.Item(119, 6).Value = $false
$ws1.Cells.Item(119, 7).Value = 25.91095506008519
$ws1.Cells.Item(119, 8).Value = 22
$ws1.Cells.Item(119, 9).Value = -26.74144946244506
$ws1.Cells.Item(119, 10).Value = 52
$ws1.Cells.Item(119, 11).Value = $true
$ws1.Cells.Item(119, 12).Value = 18
$ws1.Cells.Item(119, 13).Value = 7.894727896798699
$ws1.Cells.Item(119, 14).Value = 51.17455750077703
$ws1.Cells.Item(119, 15).Value = 33.15833033749053
$ws1.Cells.Item(119, 16).Value = -18.01622716328649

# Row 120
$ws1.Cells.Item(120, 1).Value = 'twelve_bar_breakout'
$ws1.Cells.Item(120, 2).Value = 1443
$ws1.Cells.Item(120, 4).Value = 1.167674422264099
$ws1.Cells.Item(120, 5).Value = 91.25876720989737
$ws1.Cells.Item(120, 6).Value = $true
$ws1.Cells.Item(120, 7).Value = 145.1483737699501
$ws1.Cells.Item(120, 8).Value = 45
$ws1.Cells.Item(120, 9).Value = -11.34868645685374
$ws1.Cells.Item(120, 10).Value = 0
$ws1.Cells.Item(120, 11).Value = $true
$ws1.Cells.Item(120, 12).Value = 50
$ws1.Cells.Item(120, 13).Value = 116.6439700901353
$ws1.Cells.Item(120, 14).Value = 53.8896065600527
$ws1.Cells.Item(120, 15).Value = 25.38520288023793
$ws1.Cells.Item(120, 16).Value = -28.50440367981477

# Row 121
$ws1.Cells.Item(121, 1).Value = 'twelve_bar_breakout'
$ws1.Cells.Item(121, 2).Value = 1586
$ws1.Cells.Item(121, 4).Value = 2.045235872268677
$ws1.Cells.Item(121, 5).Value = 170.7654545325256
$ws1.Cells.Item(121, 6).Value = $true
$ws1.Cells.Item(121, 7).Value = 180.8875514218833
$ws1.Cells.Item(121, 8).Value = 49
$ws1.Cells.Item(121, 9).Value = -5.131032637375474
$ws1.Cells.Item(121, 10).Value = 0
$ws1.Cells.Item(121, 11).Value = $false
$ws1.Cells.Item(121, 14).Value = 10.12209688935769

# Row 122
$ws1.Cells.Item(122, 1).Value = 'twelve_bar_breakout'
$ws1.Cells.Item(122, 2).Value = 1629
$ws1.Cells.Item(122, 4).Value = 4.958136081695557
$ws1.Cells.Item(122, 5).Value = 97.29877913610977
$ws1.Cells.Item(122, 6).Value = $true
$ws1.Cells.Item(122, 7).Value = 121.3799745153539
$ws1.Cells.Item(122, 8).Value = 44
$ws1.Cells.Item(122, 9).Value = -10.64718623138125
$ws1.Cells.Item(122, 10).Value = 1
$ws1.Cells.Item(122, 11).Value = $true
$ws1.Cells.Item(122, 12).Value = 16
$ws1.Cells.Item(122, 13).Value = 5.219194332737223
$ws1.Cells.Item(122, 14).Value = 24.08119537924409
$ws1.Cells.Item(122, 15).Value = -92.07958480337255
$ws1.Cells.Item(122, 16).Value = -116.1607801826166

# Row 123
$ws1.Cells.Item(123, 1).Value = 'twelve_bar_breakout'
$ws1.Cells.Item(123, 2).Value = 1659
$ws1.Cells.Item(123, 4).Value = 8.138830184936523
$ws1.Cells.Item(123, 5).Value = 77.64977358111182
$ws1.Cells.Item(123, 6).Value = $true
$ws1.Cells.Item(123, 7).Value = 91.95766133250339
$ws1.Cells.Item(123, 8).Value = 49
$ws1.Cells.Item(123, 9).Value = -6.740046708272708
$ws1.Cells.Item(123, 10).Value = 3
$ws1.Cells.Item(123, 11).Value = $true
$ws1.Cells.Item(123, 12).Value = 40
$ws1.Cells.Item(123, 13).Value = 64.4246865135991
$ws1.Cells.Item(123, 14).Value = 14.30788775139158
$ws1.Cells.Item(123, 15).Value = -13.22508706751272
$ws1.Cells.Item(123, 16).Value = -27.53297481890429

# Row 124
$ws1.Cells.Item(124, 1).Value = 'twelve_bar_breakout'
$ws1.Cells.Item(124, 2).Value = 1687
$ws1.Cells.Item(124, 4).Value = 13.54159927368164
$ws1.Cells.Item(124, 5).Value = -1.841388182990319
$ws1.Cells.Item(124, 6).Value = $false
$ws1.Cells.Item(124, 7).Value = 23.55225290603096
$ws1.Cells.Item(124, 8).Value = 32
$ws1.Cells.Item(124, 9).Value = -28.11356173606627
$ws1.Cells.Item(124, 10).Value = 46
$ws1.Cells.Item(124, 11).Value = $true
$ws1.Cells.Item(124, 12).Value = 12
$ws1.Cells.Item(124, 13).Value = -1.176768363950115
$ws1.Cells.Item(124, 14).Value = 25.39364108902128
$ws1.Cells.Item(124, 15).Value = 0.6646198190402033
$ws1.Cells.Item(124, 16).Value = -24.72902126998108

# Row 125
$ws1.Cells.Item(125, 1).Value = 'twelve_bar_breakout'
$ws1.Cells.Item(125, 2).Value = 1717
$ws1.Cells.Item(125, 4).Value = 16.22738838195801
$ws1.Cells.Item(125, 5).Value = -26.10772435963363
$ws1.Cells.Item(125, 6).Value = $false
$ws1.Cells.Item(125, 7).Value = 3.103164775067916
$ws1.Cells.Item(125, 8).Value = 2
$ws1.Cells.Item(125, 9).Value = -40.01145980675746
$ws1.Cells.Item(125, 10).Value = 16
$ws1.Cells.Item(125, 11).Value = $true
$ws1.Cells.Item(125, 12).Value = 9
$ws1.Cells.Item(125, 13).Value = -34.09890894558833
$ws1.Cells.Item(125, 14).Value = 29.21088913470155
$ws1.Cells.Item(125, 15).Value = -7.991184585954695
$ws1.Cells.Item(125, 16).Value = -37.20207372065624

# Row 126
$ws1.Cells.Item(126, 1).Value = 'twelve_bar_breakout'
$ws1.Cells.Item(126, 2).Value = 1797
$ws1.Cells.Item(126, 4).Value = 17.76359748840332
$ws1.Cells.Item(126, 5).Value = 40.17193370170494
$ws1.Cells.Item(126, 6).Value = $true
$ws1.Cells.Item(126, 7).Value = 42.63958321080995
$ws1.Cells.Item(126, 8).Value = 49
$ws1.Cells.Item(126, 9).Value = -13.01078632773804
$ws1.Cells.Item(126, 10).Value = 18
$ws1.Cells.Item(126, 11).Value = $true
$ws1.Cells.Item(126, 12).Value = 3
$ws1.Cells.Item(126, 13).Value = -1.022672987846014
$ws1.Cells.Item(126, 14).Value = 2.467649509105016
$ws1.Cells.Item(126, 15).Value = -41.19460668955095
$ws1.Cells.Item(126, 16).Value = -43.66225619865597

# Row 127
$ws1.Cells.Item(127, 1).Value = 'twelve_bar_breakout'
$ws1.Cells.Item(127, 2).Value = 1837
$ws1.Cells.Item(127, 4).Value = 19.9660758972168
$ws1.Cells.Item(127, 5).Value = -8.049274864458011
$ws1.Cells.Item(127, 6).Value = $false
$ws1.Cells.Item(127, 7).Value = 28.92810605526275
$ws1.Cells.Item(127, 8).Value = 13
$ws1.Cells.Item(127, 9).Value = -25.08067941322744
$ws1.Cells.Item(127, 10).Value = 45
$ws1.Cells.Item(127, 11).Value = $true
$ws1.Cells.Item(127, 12).Value = 32
$ws1.Cells.Item(127, 13).Value = 3.950296399284122
$ws1.Cells.Item(127, 14).Value = 36.97738091972076
$ws1.Cells.Item(127, 15).Value = 11.99957126374213
$ws1.Cells.Item(127, 16).Value = -24.97780965597862

# Row 128
$ws1.Cells.Item(128, 1).Value = 'twelve_bar_breakout'
$ws1.Cells.Item(128, 2).Value = 1840
$ws1.Cells.Item(128, 4).Value = 21.51316833496094
$ws1.Cells.Item(128, 5).Value = -19.8285039864897
$ws1.Cells.Item(128, 6).Value = $false
$ws1.Cells.Item(128, 7).Value = 19.65640349685234
$ws1.Cells.Item(128, 8).Value = 10
$ws1.Cells.Item(128, 9).Value = -30.46840810646527
$ws1.Cells.Item(128, 10).Value = 42
$ws1.Cells.Item(128, 11).Value = $true
$ws1.Cells.Item(128, 12).Value = 29
$ws1.Cells.Item(128, 13).Value = -3.525158399219228
$ws1.Cells.Item(128, 14).Value = 39.48490748334204
$ws1.Cells.Item(128, 15).Value = 16.30334558727047
$ws1.Cells.Item(128, 16).Value = -23.18156189607157

# Update Summary sheet row 4 (twelve_bar_breakout) with corrected aggregate stats
$ws2 = $wb.Worksheets.Item(2)
$ws2.Cells.Item(4, 2).Value = 28
$ws2.Cells.Item(4, 3).Value = 60.71428571428571
$ws2.Cells.Item(4, 4).Value = 87.71201210716326
$ws2.Cells.Item(4, 5).Value = 37.28946963091894
$ws2.Cells.Item(4, 6).Value = 148.0087538007874
$ws2.Cells.Item(4, 7).Value = 156.8062516894483
$ws2.Cells.Item(4, 8).Value = -28.42725244146962
$ws2.Cells.Item(4, 9).Value = 69.09423958228508
$ws2.Cells.Item(4, 10).Value = 71.42857142857143
$ws2.Cells.Item(4, 11).Value = 24.1
$ws2.Cells.Item(4, 12).Value = 50
$ws2.Cells.Item(4, 13).Value = -15.17611900260404
